$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Baza podataka" to "Analiza sadrzaja" in cell D3
$ws.Range("D3").Value = "Analiza sadržaja"

# Update the selection to reflect the edited cell
$ws.Range("D3").Select()
